$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "GEN2810"
$ws.Range("C2").Value = "CIE2802"
$ws.Range("D2").Value = "CIE4818"

$ws.Range("B3").Value = "POW1804"
$ws.Range("C3").Value = "GEN0809"
$ws.Range("D3").Value = "GEN1805"

$ws.Range("B4").Value = "GEN0807"
$ws.Range("C4").Value = "CIE1808"
$ws.Range("D4").Value = "MEC0811"

$ws.Range("B5").Value = "CIE1803"
$ws.Range("C5").Value = "GEN1809"
$ws.Range("D5").Value = "GEN0802"

$ws.Range("B6").Value = "GEN1801"
$ws.Range("C6").Value = "CIE3804"
$ws.Range("D6").Value = "GEN0806"

$ws.Range("B7").Value = "-"
$ws.Range("C7").Value = "GEN0810"
$ws.Range("D7").Value = "GEN0801"
